$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Input")
$ws.Cells.Item(2, 4).Value = 0.291
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0.1387570603255635
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(3, 4).Value = 0.291
$ws.Cells.Item(3, 5).Value = 0.439
$ws.Cells.Item(3, 6).Value = 0.369
$ws.Cells.Item(3, 8).Value = 0.06937853016278178
$ws.Cells.Item(4, 5).Value = 2.838766910623198
$ws.Cells.Item(4, 6).Value = 2.129075182967398
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 8).Value = 0.7096917276557995
$ws.Cells.Item(6, 4).Value = 192916.6738549898
$ws.Cells.Item(6, 5).Value = 210970.1272594048
$ws.Cells.Item(6, 6).Value = 166736.5233268229
$ws.Cells.Item(6, 7).Value = 242931.3367324624
$ws.Cells.Item(6, 8).Value = 214673.3958739777
$ws.Cells.Item(6, 9).Value = 112413.8842600082
$ws.Cells.Item(6, 10).Value = 20650.28433098331
$ws.Cells.Item(6, 11).Value = 13487.53212053503
$ws.Cells.Item(6, 12).Value = 2818.130988676962
$ws.Cells.Item(7, 4).Value = 14780.70269999836
$ws.Cells.Item(7, 5).Value = 16944.47835473503
$ws.Cells.Item(7, 6).Value = 22345.38578261545
$ws.Cells.Item(7, 7).Value = 17199.51802080974
$ws.Cells.Item(7, 8).Value = 26165.95440801004
$ws.Cells.Item(7, 9).Value = 4353.302492789441
$ws.Cells.Item(7, 10).Value = 2997.621919013707
$ws.Cells.Item(7, 11).Value = 3152.929586618577
$ws.Cells.Item(7, 12).Value = 1690.878593206179
$ws.Cells.Item(8, 5).Value = 1064.831612009241
$ws.Cells.Item(8, 6).Value = 6725.252286374135
$ws.Cells.Item(8, 9).Value = 8876.563756593521
$ws.Cells.Item(8, 10).Value = 57092.07505593903
$ws.Cells.Item(8, 11).Value = 70658.16024421096
$ws.Cells.Item(8, 12).Value = 50726.35779618536
$ws.Cells.Item(9, 6).Value = 17006.0365338585
$ws.Cells.Item(9, 10).Value = 19180.74932785417
$ws.Cells.Item(9, 11).Value = 17799.03093592133
$ws.Cells.Item(9, 12).Value = 1127.252395470788
$ws.Cells.Item(10, 4).Value = 247.2053541453888
$ws.Cells.Item(10, 5).Value = 252.4680876084153
$ws.Cells.Item(10, 6).Value = 236.8561238959894
$ws.Cells.Item(10, 7).Value = 216.2603955202925
$ws.Cells.Item(10, 8).Value = 603.826992792788
$ws.Cells.Item(10, 9).Value = 917.3119876812468
$ws.Cells.Item(10, 10).Value = 780.2207363646929
$ws.Cells.Item(10, 11).Value = 897.5870941873752
$ws.Cells.Item(10, 12).Value = 865.7888843438553
$ws.Cells.Item(11, 4).Value = 2.959138178052658
$ws.Cells.Item(11, 5).Value = 4.644969634132857
$ws.Cells.Item(11, 6).Value = 12.99209315116792
$ws.Cells.Item(11, 7).Value = 47.43384644502824
$ws.Cells.Item(11, 8).Value = 188.0816000012074
$ws.Cells.Item(11, 9).Value = 434.3660050120391
$ws.Cells.Item(11, 10).Value = 464.2004285862787
$ws.Cells.Item(11, 11).Value = 696.292634828991
$ws.Cells.Item(11, 12).Value = 865.7888843438554
$ws.Cells.Item(12, 4).Value = 74741.53624898972
$ws.Cells.Item(12, 5).Value = 61982.42202692643
$ws.Cells.Item(12, 6).Value = 55542.63651569869
$ws.Cells.Item(12, 7).Value = 42540.20368103011
$ws.Cells.Item(12, 8).Value = 13457.30601074946
$ws.Cells.Item(12, 9).Value = 12774.14483266342
$ws.Cells.Item(12, 10).Value = 460.9610785447496
$ws.Cells.Item(12, 11).Value = 151.9050674479426
$ws.Cells.Item(12, 12).Value = 218.347628920442
$ws.Cells.Item(13, 4).Value = 14095.62103188714
$ws.Cells.Item(13, 5).Value = 11850.72723030543
$ws.Cells.Item(13, 6).Value = 10619.47264199522
$ws.Cells.Item(13, 7).Value = 8133.472905052301
$ws.Cells.Item(13, 8).Value = 2572.969199539521
$ws.Cells.Item(13, 9).Value = 2442.352219578415
$ws.Cells.Item(13, 10).Value = 88.13343891673202
$ws.Cells.Item(13, 11).Value = 29.04348459382046
$ws.Cells.Item(13, 12).Value = 41.74696804516627
$ws.Cells.Item(14, 4).Value = 3897.191613117493
$ws.Cells.Item(14, 5).Value = 4132.161468461761
$ws.Cells.Item(14, 6).Value = 3702.842434379913
$ws.Cells.Item(14, 7).Value = 2836.013578735341
$ws.Cells.Item(14, 8).Value = 897.153734049964
$ws.Cells.Item(14, 9).Value = 851.6096555108948
$ws.Cells.Item(14, 10).Value = 30.73073856964998
$ws.Cells.Item(14, 11).Value = 10.1270044965295
$ws.Cells.Item(14, 12).Value = 14.55650859469613
$ws.Cells.Item(15, 4).Value = 910.5008189836789
$ws.Cells.Item(15, 5).Value = 757.64294283116
$ws.Cells.Item(15, 6).Value = 1663.738980971423
$ws.Cells.Item(15, 7).Value = 364.6705343690258
$ws.Cells.Item(15, 8).Value = 182.1001637967357
$ws.Cells.Item(15, 9).Value = 18.19838969253021
$ws.Cells.Item(15, 10).Value = 107.0816917673511
$ws.Cells.Item(15, 11).Value = 108.4417981664033
$ws.Cells.Item(16, 4).Value = 215.7244202622394
$ws.Cells.Item(16, 5).Value = 67.70714091025418
$ws.Cells.Item(16, 6).Value = 8.698905109489051
$ws.Cells.Item(16, 7).Value = 15.03369881730334
$ws.Cells.Item(17, 4).Value = 51.2022008834085
$ws.Cells.Item(17, 5).Value = 136.0671654928556
$ws.Cells.Item(17, 6).Value = 127.2946502852885
$ws.Cells.Item(17, 7).Value = 103.0385147360157
$ws.Cells.Item(17, 8).Value = 434.9257766255786
$ws.Cells.Item(17, 9).Value = 3481.229639198534
$ws.Cells.Item(17, 10).Value = 2424.068808726991
$ws.Cells.Item(17, 11).Value = 53.74897660246722
$ws.Cells.Item(18, 4).Value = 2.694852678074132
$ws.Cells.Item(18, 5).Value = 19.43816649897937
$ws.Cells.Item(18, 6).Value = 31.82366257132213
$ws.Cells.Item(18, 7).Value = 52.97610012134483
$ws.Cells.Item(18, 8).Value = 400.5895311025066
$ws.Cells.Item(18, 9).Value = 1202.450784056603
$ws.Cells.Item(18, 10).Value = 2041.31149278085
$ws.Cells.Item(18, 11).Value = 570.3973027200665
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 2.913685506673966
$ws.Cells.Item(19, 8).Value = 22.89083034871465
$ws.Cells.Item(19, 9).Value = 169.8744194963003
$ws.Cells.Item(19, 10).Value = 202.2646129602856
$ws.Cells.Item(19, 11).Value = 34.00445458523485
$ws.Cells.Item(20, 6).Value = 0.1514277794960627
$ws.Cells.Item(20, 7).Value = 0.9404135511086561
$ws.Cells.Item(20, 8).Value = 5.845005522117619
$ws.Cells.Item(20, 9).Value = 36.39677938506041
$ws.Cells.Item(20, 10).Value = 214.1633835347021
$ws.Cells.Item(22, 6).Value = 0.1642383725553879
$ws.Cells.Item(22, 9).Value = 4.372857261038206
$ws.Cells.Item(22, 10).Value = 26.23149074613449
$ws.Cells.Item(22, 11).Value = 26.2134245251534
$ws.Cells.Item(22, 12).Value = 23.99259870269209
$ws.Cells.Item(23, 7).Value = 0.1642383725553879
$ws.Cells.Item(23, 9).Value = 35.38039056658186
$ws.Cells.Item(23, 10).Value = 212.2366069459972
$ws.Cells.Item(23, 11).Value = 212.0904347944229
$ws.Cells.Item(23, 12).Value = 194.1219349581451
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 8).Value = 0.02006967768849931
$ws.Cells.Item(24, 9).Value = 0.09413062351126328
$ws.Cells.Item(24, 10).Value = 0.09413062351126328
$ws.Cells.Item(24, 11).Value = 4.485953227149231
$ws.Cells.Item(24, 12).Value = 4.391822603637968
$ws.Cells.Item(25, 7).Value = 0.8557329410114844
$ws.Cells.Item(25, 8).Value = 0.1623819376614945
$ws.Cells.Item(25, 9).Value = 0.7616023175002212
$ws.Cells.Item(25, 10).Value = 0.7616023175002212
$ws.Cells.Item(25, 11).Value = 36.29543974693469
$ws.Cells.Item(25, 12).Value = 35.53383742943447
$ws.Cells.Item(26, 8).Value = 0.1514277794960627
$ws.Cells.Item(26, 9).Value = 0.9404135511086561
$ws.Cells.Item(26, 10).Value = 5.845005522117619
$ws.Cells.Item(26, 11).Value = 36.39677938506041

$ws = $wb.Worksheets.Item("Output")
$ws.Cells.Item(2, 4).Value = 0.582
$ws.Cells.Item(2, 6).Value = 0.369
$ws.Cells.Item(2, 7).Value = 0.1387570603255635
$ws.Cells.Item(2, 8).Value = 0.06937853016278178
$ws.Cells.Item(4, 4).Value = 8842.302975664505
$ws.Cells.Item(4, 5).Value = 9745.513762208007
$ws.Cells.Item(4, 6).Value = 9040.626485100391
$ws.Cells.Item(4, 7).Value = 11890.34801039957
$ws.Cells.Item(4, 8).Value = 11055.61496935538
$ws.Cells.Item(4, 9).Value = 6241.956654381645
$ws.Cells.Item(4, 10).Value = 4354.878510122691
$ws.Cells.Item(4, 11).Value = 4588.738687813712
$ws.Cells.Item(4, 12).Value = 2465.300988894608
$ws.Cells.Item(5, 4).Value = 225.2959999999999
$ws.Cells.Item(5, 7).Value = 239.6965100910401
$ws.Cells.Item(5, 8).Value = 722.1218135146563
$ws.Cells.Item(5, 9).Value = 1238.228493674559
$ws.Cells.Item(5, 10).Value = 1143.189069885189
$ws.Cells.Item(5, 11).Value = 1469.306387856179
$ws.Cells.Item(5, 12).Value = 1601.709436036133
$ws.Cells.Item(6, 5).Value = 3203.4386870973
$ws.Cells.Item(6, 6).Value = 2870.611131015129
$ws.Cells.Item(6, 7).Value = 2198.606149492032
$ws.Cells.Item(6, 8).Value = 695.5142004649985
$ws.Cells.Item(6, 9).Value = 660.2063684081443
$ws.Cells.Item(6, 10).Value = 23.82386011980714
$ws.Cells.Item(6, 11).Value = 7.850912467045363
$ws.Cells.Item(6, 12).Value = 11.28486462526178
$ws.Cells.Item(7, 4).Value = 910.5008189836789
$ws.Cells.Item(7, 5).Value = 757.64294283116
$ws.Cells.Item(7, 6).Value = 1663.738980971423
$ws.Cells.Item(7, 7).Value = 364.6705343690258
$ws.Cells.Item(7, 8).Value = 182.1001637967357
$ws.Cells.Item(7, 9).Value = 18.19838969253021
$ws.Cells.Item(7, 10).Value = 107.0816917673511
$ws.Cells.Item(7, 11).Value = 108.4417981664033
$ws.Cells.Item(8, 4).Value = 215.7244202622394
$ws.Cells.Item(8, 5).Value = 67.70714091025418
$ws.Cells.Item(8, 6).Value = 8.698905109489051
$ws.Cells.Item(8, 7).Value = 15.03369881730334
$ws.Cells.Item(9, 4).Value = 53.48897587023141
$ws.Cells.Item(9, 5).Value = 154.325
$ws.Cells.Item(9, 6).Value = 157.929572220504
$ws.Cells.Item(9, 7).Value = 157.929572220504
$ws.Cells.Item(9, 8).Value = 853.3796387440843
$ws.Cells.Item(9, 9).Value = 4827.211827159916
$ws.Cells.Item(9, 10).Value = 4638.06210498095
$ws.Cells.Item(9, 11).Value = 654.7105034631414
$ws.Cells.Item(10, 6).Value = 0.1514277794960627
$ws.Cells.Item(10, 7).Value = 0.9404135511086561
$ws.Cells.Item(10, 8).Value = 5.845005522117619
$ws.Cells.Item(10, 9).Value = 36.39677938506041
$ws.Cells.Item(10, 10).Value = 214.1633835347021
$ws.Cells.Item(12, 6).Value = 0.07571388974803384
$ws.Cells.Item(12, 7).Value = 0.07571388974803384
$ws.Cells.Item(12, 9).Value = 18.32624724853285
$ws.Cells.Item(12, 10).Value = 109.9337930360727
$ws.Cells.Item(12, 11).Value = 109.8580791463247
$ws.Cells.Item(12, 12).Value = 100.5508000176459
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0.3944928858062944
$ws.Cells.Item(13, 8).Value = 0.08411019467634716
$ws.Cells.Item(13, 9).Value = 0.3944928858062943
$ws.Cells.Item(13, 10).Value = 0.3944928858062943
$ws.Cells.Item(13, 11).Value = 18.80022216105269
$ws.Cells.Item(13, 12).Value = 18.4057292752464
$ws.Cells.Item(14, 8).Value = 0.1514277794960627
$ws.Cells.Item(14, 9).Value = 0.9404135511086561
$ws.Cells.Item(14, 10).Value = 5.845005522117619
$ws.Cells.Item(14, 11).Value = 36.39677938506041
$ws.Cells.Item(15, 8).Value = 0.1448054107355396
$ws.Cells.Item(15, 9).Value = 0.8992996783870737
$ws.Cells.Item(15, 10).Value = 0.8992996783870743
$ws.Cells.Item(15, 11).Value = 0.8992996783870748
$ws.Cells.Item(15, 12).Value = 171.3325081900247
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 4.703958254495745
$ws.Cells.Item(16, 11).Value = 33.99284063270574
$ws.Cells.Item(16, 12).Value = 33.99284063270575
$ws.Cells.Item(17, 7).Value = 0.3868989205103038
$ws.Cells.Item(17, 8).Value = 2.292722269632621
$ws.Cells.Item(17, 9).Value = 21.9462113051214
$ws.Cells.Item(17, 10).Value = 1401.067707571761
$ws.Cells.Item(17, 11).Value = 2404.387869716356
$ws.Cells.Item(17, 12).Value = 3336.969172259398

